$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.853.19'
$ws.Range("E2").Value = '  +0.51%  '

$ws.Range("D3").Value = '2.394.49'
$ws.Range("E3").Value = '  +0.77%  '

$ws.Range("D5").Value = '504.42'
$ws.Range("E5").Value = '  -0.99%  '

$ws.Range("D6").Value = '132.63'
$ws.Range("E6").Value = '  +2.08%  '

$ws.Range("E7").Value = '  +0.13%  '

$ws.Range("D8").Value = '0.551'
$ws.Range("E8").Value = '  -0.69%  '

$ws.Range("D9").Value = '2.402.29'
$ws.Range("E9").Value = '  +0.16%  '

$ws.Range("D10").Value = '0.0975'
$ws.Range("E10").Value = '  +0.98%  '

$ws.Range("E11").Value = '  -0.75%  '

$ws.Range("E12").Value = '  +0.31%  '

$ws.Range("E13").Value = '  -1.14%  '

$ws.Range("D14").Value = '2.819.27'
$ws.Range("E14").Value = '  +0.76%  '

$ws.Range("D15").Value = '56.798.38'
$ws.Range("E15").Value = '  +0.55%  '

$ws.Range("D16").Value = '21.70'
$ws.Range("E16").Value = '  +0.32%  '

$ws.Range("D17").Value = '0.0000134'
$ws.Range("E17").Value = '  +1.73%  '

$ws.Range("D18").Value = '2.394.89'
$ws.Range("E18").Value = '  +0.12%  '

$ws.Range("D19").Value = '10.21'
$ws.Range("E19").Value = '  -0.25%  '

$ws.Range("D20").Value = '4.06'
$ws.Range("E20").Value = '  -0.03%  '

$ws.Range("D21").Value = '309.36'
$ws.Range("E21").Value = '  -1.16%  '

$ws.Range("E22").Value = '  -0.36%  '

$ws.Range("E23").Value = '  -0.04%  '

$ws.Range("E24").Value = '  -4.58%  '

$ws.Range("E25").Value = '  +2.68%  '

$ws.Range("E26").Value = '  -0.05%  '

$ws.Range("E27").Value = '  -1.35%  '

$ws.Range("E28").Value = '  -0.96%  '

$ws.Range("D29").Value = '7.40'
$ws.Range("E29").Value = '  +2.56%  '

$ws.Range("D30").Value = '175.53'
$ws.Range("E30").Value = '  +0.39%  '

$ws.Range("E31").Value = '  +0.76%  '

$ws.Range("E32").Value = '  -1.64%  '

$ws.Range("E33").Value = '  +0.38%  '

$ws.Range("D34").Value = '5.86'
$ws.Range("E34").Value = '  -4.61%  '

$ws.Range("E35").Value = '  +0.18%  '

$ws.Range("E36").Value = '  +0.28%  '

$ws.Range("D37").Value = '17.90'
$ws.Range("E37").Value = '  +0.78%  '

$ws.Range("D38").Value = '1.19'
$ws.Range("E38").Value = '  -2.05%  '

$ws.Range("D39").Value = '3.81'
$ws.Range("E39").Value = '  +1.54%  '

$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").Value = '36.85'
$ws.Range("E40").Value = '  +2.78%  '

$ws.Range("B41").Value = 'SuiNetwork'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D41").Value = '0.824'
$ws.Range("E41").Value = '  +4.61%  '

$ws.Range("E42").Value = '  +0.30%  '

$ws.Range("D43").Value = '131.15'
$ws.Range("E43").Value = '  -2.30%  '

$ws.Range("D44").Value = '3.36'
$ws.Range("E44").Value = '  +0.20%  '

$ws.Range("D45").Value = '4.85'
$ws.Range("E45").Value = '  -1.15%  '

$ws.Range("D47").Value = '0.0913'
$ws.Range("E47").Value = '  +1.36%  '

$ws.Range("D48").Value = '250.39'
$ws.Range("E48").Value = '  -2.13%  '

$ws.Range("E50").Value = '  +1.05%  '

$ws.Range("D51").Value = '17.03'
$ws.Range("E51").Value = '  +7.13%  '
